$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DesignMatrix")

# Row 6 is the "Bulk=1" row. Columns BY:DC are the "Single-cell" block.
# Clear those cells in the Bulk=1 row to 0 (they were 1).
$ws.Range("BY6:DC6").Value = 0

# Reflect the resulting view/selection state (scrolled right, new active cell).
$ws.Activate()
$excel.Goto($ws.Range("BH1"), $true)
$ws.Range("CJ16").Select()
